# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the Binance conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value2
$newText = $oldText.Replace("1000 Bs = 7.98 = 32268.13 pesos", "1000 Bs = 7.58 = 30666.59 pesos")
$newText = $newText.Replace("32268.13 pesos = 7.94 = 944.99 Bs", "30666.59 pesos = 7.56 = 945.26 Bs")
$wsHoja1.Range("A1").Value2 = $newText

# --- Sheet "tasas": update N10, O10, N12, O12 values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 132
$wsTasas.Range("O10").Value = 4047.99
$wsTasas.Range("N12").Value = 4056.95
$wsTasas.Range("O12").Value = 125.05
